# reverts LDV VTQaZ so PHEVs don't qualify from 2035
# The "VTQaZ-LDVs" sheet's "plugin hybrid vehicle" (PHEV) row (row 6) is
# updated so that years 2035-2050 (columns Q:AF) no longer qualify as ZEVs
# (value flips from 1 to 0), while 2020-2034 (columns B:P) keep qualifying.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VTQaZ-LDVs")

# Flip the plugin hybrid vehicle (row 6) qualification flag to 0 for every
# year from 2035 (column Q) through 2050 (column AF).
$ws.Range("Q6:AF6").Value2 = 0

# Reflect the editor's on-screen state: LDVs tab made active, scrolled right
# so column K is the left-most visible column, with Q6:AF6 selected.
$ws.Activate()
$ws.Range("Q6:AF6").Select()
